# Applies the edits described by the commit:
#  - bump the footer "today" date placeholder text on the slide master and
#    every slide layout from 3/11/2024 to 3/12/2024
#  - tweak wording on a few bullet / body text boxes across slides 1 and 2
#  - shrink the "Determining and filtering..." text box to its new autofit
#    height now that a word was removed from its text

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Footer date placeholder: slide master + all custom (slide) layouts
# ---------------------------------------------------------------------
$oldDate = "3/11/2024"
$newDate = "3/12/2024"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $shp = $layout.Shapes.Item($si)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2. Slide 1 body text edits
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

# "TextBox 8" - drop the word "out" and shrink the box to its new height
$tb8 = $slide1.Shapes.Item(6)
$tb8.TextFrame.TextRange.Text = "Determining and filtering aspects of the dataset which hold little value to the problem at hand, handling missing data and understanding modern biological terms and nomenclature."
$tb8.Height = 96.9375

# "TextBox 11" - "meaning" -> "semantics"
$tb11 = $slide1.Shapes.Item(8)
$tb11.TextFrame.TextRange.Text = "Finding a method to embed amino acids and V/J genes that balances the preservation of biological semantics while enabling computational efficiency."

# ---------------------------------------------------------------------
# 3. Slide 2 body text edits
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)

# "Relied on the effective..." - "meanings" -> "semantics"
$tb3a = $slide2.Shapes.Item(5)
$tb3a.TextFrame.TextRange.Text = "Relied on the effective encoding and embedding of a large quantity of string-type data and incorporated BLOSUM matrices to preserve biological semantics of amino acid sequences."

# "The resulting matrices..." - "they square matrices" -> "they are square matrices"
$tb3b = $slide2.Shapes.Item(6)
$tb3b.TextFrame.TextRange.Text = "The resulting matrices are massive in size as they are square matrices with tens of thousands of rows/columns, occupying gigabytes of storage space even when compressed."
